$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.03186166666666667"
$ws.Range("H2").Value = [double]"0.095585"
$ws.Range("I2").Value = [double]"0.0002078156820111728"
$ws.Range("J2").Value = [double]"0.0002078156820111728"
$ws.Range("M2").Value = [double]"0.1516543333333333"
$ws.Range("N2").Value = [double]"0.454963"
$ws.Range("O2").Value = [double]"0.0007044400935133411"
$ws.Range("P2").Value = [double]"0.0007044400935133412"
$ws.Range("Q2").Value = [double]"0.004831959817222222"
$ws.Range("R2").Value = [double]"0.043487638355"
$ws.Range("S2").Value = [double]"1.463936984694893E-07"
$ws.Range("T2").Value = [double]"1.463936984694893E-07"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.03186166666666667"
$ws.Range("H3").Value = [double]"0.095585"
$ws.Range("I3").Value = [double]"0.0002078156820111728"
$ws.Range("J3").Value = [double]"0.0002078156820111728"
$ws.Range("O3").Value = [double]"0.0008677905854558892"
$ws.Range("P3").Value = [double]"0.0008677905854558892"
$ws.Range("Q3").Value = [double]"0.005952428428333333"
$ws.Range("R3").Value = [double]"0.053571855855"
$ws.Range("S3").Value = [double]"1.803404923593906E-07"
$ws.Range("T3").Value = [double]"1.803404923593905E-07"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.03186166666666667"
$ws.Range("H4").Value = [double]"0.095585"
$ws.Range("I4").Value = [double]"0.0002078156820111728"
$ws.Range("J4").Value = [double]"0.0002078156820111728"
$ws.Range("K4").Value = [double]"1"
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.003907"
$ws.Range("N4").Value = [double]"0.011721"
$ws.Range("O4").Value = [double]"1.814816223752233E-05"
$ws.Range("P4").Value = [double]"1.814816223752233E-05"
$ws.Range("Q4").Value = [double]"0.0001244835316666667"
$ws.Range("R4").Value = [double]"0.001120351785"
$ws.Range("S4").Value = [double]"3.771472712640115E-09"
$ws.Range("T4").Value = [double]"3.771472712640115E-09"
$ws.Range("E5").Value = [double]"2"
$ws.Range("F5").Value = [double]"0.6666666666666666"
$ws.Range("G5").Value = [double]"0.03186166666666667"
$ws.Range("H5").Value = [double]"0.095585"
$ws.Range("I5").Value = [double]"0.0002078156820111728"
$ws.Range("J5").Value = [double]"0.0002078156820111728"
$ws.Range("M5").Value = [double]"214.941124"
$ws.Range("N5").Value = [double]"644.8233720000001"
$ws.Range("O5").Value = [double]"0.9984096211587932"
$ws.Range("P5").Value = [double]"0.9984096211587933"
$ws.Range("Q5").Value = [double]"6.848382445846668"
$ws.Range("R5").Value = [double]"61.63544201262001"
$ws.Range("S5").Value = [double]"0.0002074851763476313"
$ws.Range("T5").Value = [double]"0.0002074851763476313"
$ws.Range("I6").Value = [double]"0.0004817036408055181"
$ws.Range("J6").Value = [double]"0.0004817036408055181"
$ws.Range("M6").Value = [double]"0.1516543333333333"
$ws.Range("N6").Value = [double]"0.454963"
$ws.Range("O6").Value = [double]"0.0007044400935133411"
$ws.Range("P6").Value = [double]"0.0007044400935133412"
$ws.Range("Q6").Value = [double]"0.01120017803111111"
$ws.Range("R6").Value = [double]"0.10080160228"
$ws.Range("S6").Value = [double]"3.393313577747561E-07"
$ws.Range("T6").Value = [double]"3.393313577747561E-07"
$ws.Range("I7").Value = [double]"0.0004817036408055181"
$ws.Range("J7").Value = [double]"0.0004817036408055181"
$ws.Range("O7").Value = [double]"0.0008677905854558892"
$ws.Range("P7").Value = [double]"0.0008677905854558892"
$ws.Range("S7").Value = [double]"4.180178844708539E-07"
$ws.Range("T7").Value = [double]"4.180178844708539E-07"
$ws.Range("I8").Value = [double]"0.0004817036408055181"
$ws.Range("J8").Value = [double]"0.0004817036408055181"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.003907"
$ws.Range("N8").Value = [double]"0.011721"
$ws.Range("O8").Value = [double]"1.814816223752233E-05"
$ws.Range("P8").Value = [double]"1.814816223752233E-05"
$ws.Range("Q8").Value = [double]"0.0002885449733333334"
$ws.Range("R8").Value = [double]"0.00259690476"
$ws.Range("S8").Value = [double]"8.742035823743724E-09"
$ws.Range("T8").Value = [double]"8.742035823743724E-09"
$ws.Range("I9").Value = [double]"0.0004817036408055181"
$ws.Range("J9").Value = [double]"0.0004817036408055181"
$ws.Range("M9").Value = [double]"214.941124"
$ws.Range("N9").Value = [double]"644.8233720000001"
$ws.Range("O9").Value = [double]"0.9984096211587932"
$ws.Range("P9").Value = [double]"0.9984096211587933"
$ws.Range("Q9").Value = [double]"15.87411847781334"
$ws.Range("R9").Value = [double]"142.86706630032"
$ws.Range("S9").Value = [double]"0.0004809375495274487"
$ws.Range("T9").Value = [double]"0.0004809375495274487"
$ws.Range("G10").Value = [double]"0.1279203333333333"
$ws.Range("H10").Value = [double]"0.383761"
$ws.Range("I10").Value = [double]"0.0008343521885681821"
$ws.Range("J10").Value = [double]"0.000834352188568182"
$ws.Range("M10").Value = [double]"0.1516543333333333"
$ws.Range("N10").Value = [double]"0.454963"
$ws.Range("O10").Value = [double]"0.0007044400935133411"
$ws.Range("P10").Value = [double]"0.0007044400935133412"
$ws.Range("Q10").Value = [double]"0.01939967287144444"
$ws.Range("R10").Value = [double]"0.174597055843"
$ws.Range("S10").Value = [double]"5.87751133738031E-07"
$ws.Range("T10").Value = [double]"5.877511337380311E-07"
$ws.Range("G11").Value = [double]"0.1279203333333333"
$ws.Range("H11").Value = [double]"0.383761"
$ws.Range("I11").Value = [double]"0.0008343521885681821"
$ws.Range("J11").Value = [double]"0.000834352188568182"
$ws.Range("O11").Value = [double]"0.0008677905854558892"
$ws.Range("P11").Value = [double]"0.0008677905854558892"
$ws.Range("Q11").Value = [double]"0.02389820459366666"
$ws.Range("R11").Value = [double]"0.215083841343"
$ws.Range("S11").Value = [double]"7.240429741939853E-07"
$ws.Range("T11").Value = [double]"7.240429741939852E-07"
$ws.Range("G12").Value = [double]"0.1279203333333333"
$ws.Range("H12").Value = [double]"0.383761"
$ws.Range("I12").Value = [double]"0.0008343521885681821"
$ws.Range("J12").Value = [double]"0.000834352188568182"
$ws.Range("K12").Value = [double]"1"
$ws.Range("L12").Value = [double]"0.3333333333333333"
$ws.Range("M12").Value = [double]"0.003907"
$ws.Range("N12").Value = [double]"0.011721"
$ws.Range("O12").Value = [double]"1.814816223752233E-05"
$ws.Range("P12").Value = [double]"1.814816223752233E-05"
$ws.Range("Q12").Value = [double]"0.0004997847423333333"
$ws.Range("R12").Value = [double]"0.004498062681"
$ws.Range("S12").Value = [double]"1.514195888136719E-08"
$ws.Range("T12").Value = [double]"1.514195888136719E-08"
$ws.Range("G13").Value = [double]"0.1279203333333333"
$ws.Range("H13").Value = [double]"0.383761"
$ws.Range("I13").Value = [double]"0.0008343521885681821"
$ws.Range("J13").Value = [double]"0.000834352188568182"
$ws.Range("M13").Value = [double]"214.941124"
$ws.Range("N13").Value = [double]"644.8233720000001"
$ws.Range("O13").Value = [double]"0.9984096211587932"
$ws.Range("P13").Value = [double]"0.9984096211587933"
$ws.Range("Q13").Value = [double]"27.49534022912134"
$ws.Range("R13").Value = [double]"247.458062062092"
$ws.Range("S13").Value = [double]"0.0008330252525013687"
$ws.Range("T13").Value = [double]"0.0008330252525013687"
$ws.Range("G14").Value = [double]"153.0833153333333"
$ws.Range("H14").Value = [double]"459.249946"
$ws.Range("I14").Value = [double]"0.9984761284886152"
$ws.Range("J14").Value = [double]"0.9984761284886152"
$ws.Range("M14").Value = [double]"0.1516543333333333"
$ws.Range("N14").Value = [double]"0.454963"
$ws.Range("O14").Value = [double]"0.0007044400935133411"
$ws.Range("P14").Value = [double]"0.0007044400935133412"
$ws.Range("Q14").Value = [double]"23.21574813133311"
$ws.Range("R14").Value = [double]"208.941733181998"
$ws.Range("S14").Value = [double]"0.000703366617323359"
$ws.Range("T14").Value = [double]"0.0007033666173233591"
$ws.Range("G15").Value = [double]"153.0833153333333"
$ws.Range("H15").Value = [double]"459.249946"
$ws.Range("I15").Value = [double]"0.9984761284886152"
$ws.Range("J15").Value = [double]"0.9984761284886152"
$ws.Range("O15").Value = [double]"0.0008677905854558892"
$ws.Range("P15").Value = [double]"0.0008677905854558892"
$ws.Range("Q15").Value = [double]"28.59917805388866"
$ws.Range("R15").Value = [double]"257.392602484998"
$ws.Range("S15").Value = [double]"0.000866468184104865"
$ws.Range("T15").Value = [double]"0.000866468184104865"
$ws.Range("G16").Value = [double]"153.0833153333333"
$ws.Range("H16").Value = [double]"459.249946"
$ws.Range("I16").Value = [double]"0.9984761284886152"
$ws.Range("J16").Value = [double]"0.9984761284886152"
$ws.Range("K16").Value = [double]"1"
$ws.Range("L16").Value = [double]"0.3333333333333333"
$ws.Range("M16").Value = [double]"0.003907"
$ws.Range("N16").Value = [double]"0.011721"
$ws.Range("O16").Value = [double]"1.814816223752233E-05"
$ws.Range("P16").Value = [double]"1.814816223752233E-05"
$ws.Range("Q16").Value = [double]"0.5980965130073334"
$ws.Range("R16").Value = [double]"5.382868617066"
$ws.Range("S16").Value = [double]"1.812050677010458E-05"
$ws.Range("T16").Value = [double]"1.812050677010458E-05"
$ws.Range("G17").Value = [double]"153.0833153333333"
$ws.Range("H17").Value = [double]"459.249946"
$ws.Range("I17").Value = [double]"0.9984761284886152"
$ws.Range("J17").Value = [double]"0.9984761284886152"
$ws.Range("M17").Value = [double]"214.941124"
$ws.Range("N17").Value = [double]"644.8233720000001"
$ws.Range("O17").Value = [double]"0.9984096211587932"
$ws.Range("P17").Value = [double]"0.9984096211587933"
$ws.Range("Q17").Value = [double]"32903.8998633931"
$ws.Range("R17").Value = [double]"296135.0987705379"
$ws.Range("S17").Value = [double]"0.9968881731804168"
$ws.Range("T17").Value = [double]"0.9968881731804169"
